$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert new rows (top-to-bottom order so subsequent indices are correct)
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(15).Insert()

# Write data rows 2-16 (header row 1 untouched)
# Row 2
$ws.Cells.Item(2, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(2, 2).Value = 'Inkscape「Hershey Text」用svgフォント変換ツール開発(python)'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5389316'
$ws.Cells.Item(2, 7).Value = 315
$ws.Cells.Item(2, 8).Value = '🔥Python ◆ツール,開発'
# Row 3
$ws.Cells.Item(3, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(3, 2).Value = '自社開発のロジシステムをサポート及び開発できる方募集【PHP, Python, VBA etc】'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5389460'
$ws.Cells.Item(3, 7).Value = 305
$ws.Cells.Item(3, 8).Value = '🔥Python ◆開発 ○PHP'
# Row 4
$ws.Cells.Item(4, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(4, 2).Value = 'Excel・Accessベースの改修や追加、Pythonスクレイピングやデータ整形等の開発員募集'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5273634'
$ws.Cells.Item(4, 7).Value = 298
$ws.Cells.Item(4, 8).Value = '🔥Python ◆開発,スクレイピング'
# Row 5
$ws.Cells.Item(5, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(5, 2).Value = '【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5314730'
$ws.Cells.Item(5, 7).Value = 178
$ws.Cells.Item(5, 8).Value = '★bot ◆ツール'
# Row 6
$ws.Cells.Item(6, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(6, 2).Value = '【急募】屋上貸切露天風呂の空き状況確認システム開発'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5389645'
$ws.Cells.Item(6, 7).Value = 125
$ws.Cells.Item(6, 8).Value = '◆開発,システム開発'
# Row 7
$ws.Cells.Item(7, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(7, 2).Value = '【急募】Windowsサイネージシステム開発のプロフェッショナル募集'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5388877'
$ws.Cells.Item(7, 7).Value = 125
$ws.Cells.Item(7, 8).Value = '◆開発,システム開発'
# Row 8
$ws.Cells.Item(8, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(8, 2).Value = '【システム開発】LINEメッセージ自動処理・スプレッドシート連携'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5389694'
$ws.Cells.Item(8, 7).Value = 118
$ws.Cells.Item(8, 8).Value = '◆開発,システム開発'
# Row 9
$ws.Cells.Item(9, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(9, 2).Value = '【急募】管理システムの詳細設計・開発を依頼します!'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5389414'
$ws.Cells.Item(9, 7).Value = 115
$ws.Cells.Item(9, 8).Value = '◆開発 ◇管理'
# Row 10
$ws.Cells.Item(10, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(10, 2).Value = '【急募】ECサイトのインタラクティブな商品比較シュミレーションの開発'
$ws.Cells.Item(10, 3).Value = 'システム開発'
$ws.Cells.Item(10, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = '期限情報なし'
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5389306'
$ws.Cells.Item(10, 7).Value = 93
$ws.Cells.Item(10, 8).Value = '◆開発 ◇サイト'
# Row 11
$ws.Cells.Item(11, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(11, 2).Value = '【急募】Google Cloud WordPress管理画面ログイン設定'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5388922'
$ws.Cells.Item(11, 7).Value = 50
$ws.Cells.Item(11, 8).Value = '◇管理 ○WordPress'
# Row 12
$ws.Cells.Item(12, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(12, 2).Value = 'OR(operations research)にて最適化の仕組みの構築(社内常駐)'
$ws.Cells.Item(12, 3).Value = 'システム開発'
$ws.Cells.Item(12, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(12, 5).Value = '期限情報なし'
$ws.Cells.Item(12, 6).Value = 'https://www.lancers.jp/work/detail/5372984'
$ws.Cells.Item(12, 7).Value = 25
# Row 13
$ws.Cells.Item(13, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(13, 2).Value = '限定公開 PR 限定公開の仕事'
$ws.Cells.Item(13, 3).Value = 'システム開発'
$ws.Cells.Item(13, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(13, 5).Value = '期限情報なし'
$ws.Cells.Item(13, 6).Value = 'https://www.lancers.jp/work/detail/5385681'
$ws.Cells.Item(13, 7).Value = 25
# Row 14
$ws.Cells.Item(14, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(14, 2).Value = 'Google Ad Managerの設定支援とGoogleアドセンス・SSPの収益最大化支援'
$ws.Cells.Item(14, 3).Value = 'システム開発'
$ws.Cells.Item(14, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(14, 5).Value = '期限情報なし'
$ws.Cells.Item(14, 6).Value = 'https://www.lancers.jp/work/detail/5389241'
$ws.Cells.Item(14, 7).Value = 13
# Row 15
$ws.Cells.Item(15, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(15, 2).Value = 'Googleアナリティクスの設置方法の伝授'
$ws.Cells.Item(15, 3).Value = 'システム開発'
$ws.Cells.Item(15, 4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(15, 5).Value = '期限情報なし'
$ws.Cells.Item(15, 6).Value = 'https://www.lancers.jp/work/detail/5389676'
$ws.Cells.Item(15, 7).Value = 10
# Row 16
$ws.Cells.Item(16, 1).Value = '2025-09-08 18:26:54'
$ws.Cells.Item(16, 2).Value = '【至急】【継続案件】エラーで起動しなくなったエクセルマクロの修正をお願い致します。'
$ws.Cells.Item(16, 3).Value = 'システム開発'
$ws.Cells.Item(16, 4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(16, 5).Value = '期限情報なし'
$ws.Cells.Item(16, 6).Value = 'https://www.lancers.jp/work/detail/5389081'
$ws.Cells.Item(16, 7).Value = 10

# Rebuild hyperlinks for column F (rows 2-16) so relationship targets match the URL text
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 16; $r++) {
  $cell = $ws.Cells.Item($r, 6)
  $url = $cell.Text
  $ws.Hyperlinks.Add($cell, $url)
}
